# Bug fix in Kbar now-casting calculation (see commit message).
# Updates computed $I^K_{fc}$ values and downstream dependent sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fig. 4c")
$ws.Range("C2").Value = 1.319214963900779
$ws.Range("C3").Value = 1.352917872027837
$ws.Range("C4").Value = 1.822846887878908
$ws.Range("C5").Value = 1.365562930176478
$ws.Range("C6").Value = 2.284032009807097
$ws.Range("C7").Value = 1.859318468447161
$ws.Range("C8").Value = 2.50486007854281
$ws.Range("C9").Value = 1.727255000348879
$ws.Range("C10").Value = 2.153623006296188
$ws.Range("C11").Value = 0.8562377486487824
$ws.Range("C12").Value = 2.203508468714841
$ws.Range("C13").Value = 1.565967862858648
$ws.Range("C14").Value = 1.642785309239926
$ws.Range("C15").Value = 1.939856566834212
$ws.Range("C16").Value = 0.8411184920888181
$ws.Range("C17").Value = 1.221281379481724
$ws.Range("C18").Value = 2.901500037896975
$ws.Range("C19").Value = 4.247749710805369
$ws.Range("C20").Value = 1.807205341747676
$ws.Range("C21").Value = 1.858418584659465
$ws.Range("C22").Value = 2.033450405181828
$ws.Range("C23").Value = 1.730192540405501
$ws.Range("C24").Value = 2.15697541653352
$ws.Range("C25").Value = 2.118638154129988
$ws.Range("C26").Value = 1.974188063803748
$ws.Range("C27").Value = 2.050588083986683
$ws.Range("C28").Value = 2.028326688383942
$ws.Range("C29").Value = 1.585553385331696
$ws.Range("C30").Value = 2.804009740469744
$ws.Range("C31").Value = 2.453655800576705
$ws.Range("C32").Value = 3.147659891102718
$ws.Range("C33").Value = 1.895659524163687
$ws.Range("C34").Value = 1.882229673246511
$ws.Range("C35").Value = 1.866611804813323
$ws.Range("C36").Value = 2.304019038126573
$ws.Range("C37").Value = 3.068859199568302
$ws.Range("C38").Value = 5.038693393248002
$ws.Range("C39").Value = 2.563389789559973
$ws.Range("C40").Value = 3.376043313878982
$ws.Range("C41").Value = 2.418229892875246
$ws.Range("C42").Value = 2.338695644741247
$ws.Range("C43").Value = 2.234326960909044
$ws.Range("C44").Value = 2.974159260662775
$ws.Range("C45").Value = 2.783010053574825
$ws.Range("C46").Value = 2.13068690950413
$ws.Range("C47").Value = 3.672033848623827
$ws.Range("C48").Value = 2.308079875176504
$ws.Range("C49").Value = 2.528999603541364
$ws.Range("C50").Value = 4.264918499627088

$ws = $wb.Worksheets.Item("Fig. 4c regression")
$ws.Range("B2").Value = 0.02294010474866214
$ws.Range("B3").Value = 1.416327319838081
$ws.Range("B4").Value = 0.5588345907055456
$ws.Range("B5").Value = 0.0000300241400273786
$ws.Range("B6").Value = 0.004965510803936371

$ws = $wb.Worksheets.Item("World values")
$ws.Range("B3").Value = 2.039298200777883

$ws = $wb.Worksheets.Item("Fig. 4e")
$ws.Range("B2").Value = 0.09133431721496336
$ws.Range("B3").Value = 0.2287247897012488
$ws.Range("B4").Value = 0.5040986771504566
$ws.Range("B5").Value = -0.2641016538850014
$ws.Range("B6").Value = 0.05454581368534421
$ws.Range("B7").Value = 0.9045410716600322
$ws.Range("B8").Value = 0.1238102265676936
$ws.Range("B9").Value = 0.2056271466662294
$ws.Range("B10").Value = 0.00552089280691645
$ws.Range("B11").Value = -0.005394567739185516
$ws.Range("B12").Value = 0.1369876277364692
$ws.Range("B13").Value = 0.3773556681589778
$ws.Range("B14").Value = 0.0381677727006323
$ws.Range("B15").Value = -0.2516722440364716
$ws.Range("B16").Value = 0.4087001668668862
$ws.Range("B17").Value = -0.08014865892999604
$ws.Range("B18").Value = -0.08848081107889362
$ws.Range("B19").Value = 0.7337836334604819
$ws.Range("B20").Value = -0.04999169372100996
$ws.Range("B21").Value = -0.1208240874597795
$ws.Range("B22").Value = -0.4010389827595632
$ws.Range("B23").Value = 0.2152180814814271
$ws.Range("B24").Value = -0.4103420827397302
$ws.Range("B25").Value = -0.07303891759728816
$ws.Range("B26").Value = 0.3184447135015584
$ws.Range("B27").Value = 0.1849733464461891
$ws.Range("B28").Value = 0.05610124083692785
$ws.Range("B29").Value = 0.737817342205254
$ws.Range("B30").Value = 0.3526221282682845
$ws.Range("B31").Value = 0.4340535556019213
$ws.Range("B32").Value = -0.8678129270968798
$ws.Range("B33").Value = 0.3109273655912662
$ws.Range("B34").Value = 0.5881399613091615
$ws.Range("B35").Value = -0.09287982666249232
$ws.Range("B36").Value = -0.1643730320924714
$ws.Range("B37").Value = -0.8856284640960508
$ws.Range("B38").Value = -0.2162125688939541
$ws.Range("B39").Value = 0.17043009386258
$ws.Range("B40").Value = -0.03244852255075728
$ws.Range("B41").Value = -0.002871672268772771
$ws.Range("B42").Value = -0.5127051108363867
$ws.Range("B43").Value = 0.1220492764338675
$ws.Range("B44").Value = 0.04383869101526796
$ws.Range("B45").Value = -0.1122062263676609
$ws.Range("B46").Value = -0.09239572367596294
$ws.Range("B47").Value = -0.4355688940793756
$ws.Range("B48").Value = -0.1660722860819145
$ws.Range("B49").Value = 0.07744511915260491
$ws.Range("B50").Value = 0.1133365770558555

$ws = $wb.Worksheets.Item("Fig. 4b")
$ws.Range("D2").Value = 2.234326960909044
$ws.Range("D3").Value = 2.563389789559973
$ws.Range("D4").Value = 3.376043313878982
$ws.Range("D5").Value = 1.565967862858648
$ws.Range("D6").Value = 2.153623006296188
$ws.Range("D7").Value = 5.038693393248002
$ws.Range("D8").Value = 2.308079875176504
$ws.Range("D9").Value = 2.50486007854281
$ws.Range("D10").Value = 2.050588083986683
$ws.Range("D11").Value = 2.028326688383942
$ws.Range("D12").Value = 2.338695644741247
$ws.Range("D13").Value = 2.974159260662775
$ws.Range("D14").Value = 2.118638154129988
$ws.Range("D15").Value = 1.585553385331696
$ws.Range("D16").Value = 3.068859199568302
$ws.Range("D17").Value = 1.882229673246511
$ws.Range("D18").Value = 1.866611804813323
$ws.Range("D19").Value = 4.247749710805369
$ws.Range("D20").Value = 1.939856566834212
$ws.Range("D21").Value = 1.807205341747676
$ws.Range("D22").Value = 1.365562930176478
$ws.Range("D23").Value = 2.528999603541364
$ws.Range("D24").Value = 1.352917872027837
$ws.Range("D25").Value = 1.895659524163687
$ws.Range("D26").Value = 2.804009740469744
$ws.Range("D27").Value = 2.453655800576705
$ws.Range("D28").Value = 2.15697541653352
$ws.Range("D29").Value = 4.264918499627088
$ws.Range("D30").Value = 2.901500037896975
$ws.Range("D31").Value = 3.147659891102718
$ws.Range("D32").Value = 0.8562377486487824
$ws.Range("D33").Value = 2.783010053574825
$ws.Range("D34").Value = 3.672033848623827
$ws.Range("D35").Value = 1.858418584659465
$ws.Range("D36").Value = 1.730192540405501
$ws.Range("D37").Value = 0.8411184920888181
$ws.Range("D38").Value = 1.642785309239926
$ws.Range("D39").Value = 2.418229892875246
$ws.Range("D40").Value = 1.974188063803748
$ws.Range("D41").Value = 2.033450405181828
$ws.Range("D42").Value = 1.221281379481724
$ws.Range("D43").Value = 2.304019038126573
$ws.Range("D44").Value = 2.13068690950413
$ws.Range("D45").Value = 1.822846887878908
$ws.Range("D46").Value = 1.859318468447161
$ws.Range("D47").Value = 1.319214963900779
$ws.Range("D48").Value = 1.727255000348879
$ws.Range("D49").Value = 2.203508468714841
$ws.Range("D50").Value = 2.284032009807097

$ws = $wb.Worksheets.Item("elasticity E_K_fc")
$ws.Range("B2").Value = 0.2853815620769785
$ws.Range("B3").Value = -0.2241448428323133
$ws.Range("B4").Value = 0.5296384308481413
$ws.Range("B5").Value = 0.00009102840441571928
$ws.Range("B6").Value = 0.06666652223473167

